# v0.94 - Clock server changed counter is added to the ralog analysis
#
# Adds two new switch pairs to the "switch_rack" sheet:
#   - VC01 / Brocade G630 pair, inserted right above the "OST Citrix_VDI" block
#   - VC02 / Brocade G620 pair, appended at the end of the table
# Updates the hidden _FilterDatabase defined name to the new used range and
# applies a "duplicate values" conditional format on column C (switchWwn).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("switch_rack")

# --- Insert 2 rows above the "OST Citrix_VDI" block: VC01 / Brocade G630 -----
# Inserting (rather than writing past the used range) makes Excel copy down
# the formatting of the row above, matching the workbook's existing style
# pattern for this block (default font on A/D/E, alternate font on B/C).
$ws.Rows("96:97").Insert()

$ws.Range("A96").Value = "VC01"
$ws.Range("A97").Value = "VC01"

$ws.Range("B96").Value = "o3-g630-003-vc01-f1"
$ws.Range("B97").Value = "o3-g630-004-vc01-f2"

$ws.Range("D96").Value = "Brocade G630"
$ws.Range("D97").Value = "Brocade G630"

$ws.Range("E96").Value = "1H15C140 un40-41"
$ws.Range("E97").Value = "1H15C180 un40-41"

$ws.Range("C96").Value = "10:00:88:94:71:bd:43:20"
$ws.Range("C97").Value = "10:00:88:94:71:c2:fe:40"

# --- Append 2 rows at the bottom of the table: VC02 / Brocade G620 -----------
$ws.Rows("110:111").Insert()

$ws.Range("A110").Value = "VC02"
$ws.Range("A111").Value = "VC02"

$ws.Range("B110").Value = "n4-g620-003-vc02-f1"
$ws.Range("B111").Value = "n4-g620-004-vc02-f2"

$ws.Range("D110").Value = "Brocade G620"
$ws.Range("D111").Value = "Brocade G620"

$ws.Range("C110").Value = "10:00:88:94:71:45:03:00"
$ws.Range("C111").Value = "10:00:d8:1f:cc:a1:cb:80"

$ws.Range("E110").Value = "6H6H150 un38"
$ws.Range("E111").Value = "6H6H170 un38"

# --- Update the hidden AutoFilter range (_FilterDatabase) --------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "switch_rack!_FilterDatabase") {
        $n.RefersTo = "=switch_rack!`$A`$1:`$E`$103"
    }
}

# --- Conditional formatting: highlight duplicate switchWwn values (col C) ----
$fc = $ws.Range("C1:C1048576").FormatConditions
$cond = $fc.AddUniqueValues()
$cond.DupeUnique = 1
$cond.Font.Color = 393372
$cond.Interior.Color = 13551615

# --- Restore view state (active cell / scroll position) ----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 79
$ws.Range("G117").Select()
